$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.439239263534546
$ws.Range("B1").Value = 3.276188850402832
$ws.Range("C1").Value = 2.998655557632446
$ws.Range("D1").Value = 3.404665470123291
$ws.Range("E1").Value = 1.878867626190186
